# Implementación parcial de Dublin Core
# Restructure "Obras" (sheet1) and "Referentes" (sheet2) to use Dublin Core
# style column names, split the "Dimensiones" column into separate
# width/height columns, and turn the free-text dates in "Referentes" into
# real Excel dates.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Obras")
$ws2 = $wb.Worksheets.Item("Referentes")

# ---------------------------------------------------------------------
# 1) Fix a pre-existing filename typo ("ocalo-tragedia.jpg" ->
#    "zocalo-tragedia.jpg") on the "Zócalo de la tragedia" row.
# ---------------------------------------------------------------------
$ws1.Range("H4").Value = "zocalo-tragedia.jpg"

# ---------------------------------------------------------------------
# 2) Split "Dimensiones" into "Ancho cm" (width) / "Alto cm" (height).
# ---------------------------------------------------------------------
$ws1.Range("E1").Value = "Ancho cm"
$ws1.Range("F1").Value = "Alto cm"

$ws1.Range("E2").Value = 120
$ws1.Range("F2").Value = 100
$ws1.Range("E3").Value = 160
$ws1.Range("F3").Value = 45
$ws1.Range("E4").Value = 100
$ws1.Range("F4").Value = 70

# ---------------------------------------------------------------------
# 3) Rename "Periódico" -> "Publisher" on the Referentes sheet.
# ---------------------------------------------------------------------
$ws2.Range("D1").Value = "Publisher"

# ---------------------------------------------------------------------
# 4) Rename "Fecha" -> "Date" and "Título" -> "Creator" (Obras), filling
#    the Creator value.
# ---------------------------------------------------------------------
$ws1.Range("D1").Value = "Date"
$ws1.Range("B1").Value = "Creator"
$ws1.Range("B2").Value = "Beatriz González"
$ws1.Range("B3").Value = "Beatriz González"
$ws1.Range("B4").Value = "Beatriz González"

# ---------------------------------------------------------------------
# 5) Add the "Title" column (the actual artwork titles move here).
# ---------------------------------------------------------------------
$ws1.Range("C1").Value = "Title"
$ws1.Range("C2").Value = "Los Suicidas del Sisga No 1"
$ws1.Range("C3").Value = "El Paraíso"
$ws1.Range("C4").Value = "Zócalo de la tragedia"

# ---------------------------------------------------------------------
# 6) Add the "Format.medium" column (the "Técnica" values move here).
# ---------------------------------------------------------------------
$ws1.Range("G1").Value = "Format.medium"
$ws1.Range("G2").Value = "Óleo sobre lienzo"
$ws1.Range("G3").Value = "Óleo sobre lienzo"
$ws1.Range("G4").Value = "Tipografía sobre papel"

# ---------------------------------------------------------------------
# 7) Move "Archivo"/"Referentes" to columns H/I and fill in the
#    remaining untouched values so the full table matches again.
# ---------------------------------------------------------------------
$ws1.Range("A1").Value = "ID"
$ws1.Range("H1").Value = "Archivo"
$ws1.Range("H1").Font.Bold = $true
$ws1.Range("I1").Value = "Referentes"
$ws1.Range("I1").Font.Bold = $true

$ws1.Range("A2").Value = 1
$ws1.Range("D2").Value = 1965
$ws1.Range("H2").Value = "suicidas-sisga-1.jpg"
$ws1.Range("I2").Value = 1

$ws1.Range("A3").Value = 2
$ws1.Range("D3").Value = 1997
$ws1.Range("H3").Value = "el-paraiso.jpg"
$ws1.Range("I3").Value = 2.3

$ws1.Range("A4").Value = 3
$ws1.Range("D4").Value = 1983
$ws1.Range("H4").Value = "zocalo-tragedia.jpg"
$ws1.Range("I4").Value = 4

# Column widths (best effort - the headless engine quantizes to 1/6 char)
$ws1.Columns.Item(2).ColumnWidth = 14.7369791666667
$ws1.Columns.Item(3).ColumnWidth = 23.4518229166667
$ws1.Columns.Item(4).ColumnWidth = 5.30729166666667
$ws1.Columns.Item(5).ColumnWidth = 11.7369791666667
$ws1.Columns.Item(6).ColumnWidth = 11.7369791666667
$ws1.Columns.Item(7).ColumnWidth = 19.8776041666667
$ws1.Columns.Item(8).ColumnWidth = 17.5924479166667

$ws1.Range("G2").Select()

# ---------------------------------------------------------------------
# 8) Finish the "Referentes" sheet: rename remaining headers and turn
#    the free-text dates into real Excel date values.
# ---------------------------------------------------------------------
$ws2.Range("A1").Value = "ID"
$ws2.Range("B1").Value = "Title"
$ws2.Range("C1").Value = "Date"
$ws2.Range("E1").Value = "Archivo"

$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = 'Doble suicidio en "El Sisga"'
$ws2.Range("C2").Value = 23922
$ws2.Range("C2").NumberFormat = "mm-dd-yy"
$ws2.Range("D2").Value = "El Tiempo"
$ws2.Range("E2").Value = "doble-suicidio-el-tiempo.jpg"

$ws2.Range("A3").Value = 2
$ws2.Range("B3").Value = "Una indígena y su hijo murieron en persecución"
$ws2.Range("C3").Value = 35209
$ws2.Range("D3").Value = "El Tiempo"
$ws2.Range("E3").Value = "indigena-hijo-el-tiempo.jpg"

# Re-use the exact same style record for C3 as C2 (copy/paste the format
# instead of re-assigning NumberFormat, which would otherwise create a
# duplicate - but functionally identical - style entry).
$ws2.Range("C2").Copy()
$ws2.Range("C3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws2.Range("A4").Value = 3
$ws2.Range("B4").Value = "Láminas de paisajes latinoamericanos"
$ws2.Range("E4").Value = "laminas-paisajes.jpg"

$ws2.Range("A5").Value = 4
$ws2.Range("B5").Value = "Exmilitar Mata a la Esposa de su Amigo y se Suicida"
$ws2.Range("E5").Value = "exmilitar-mata-esposa.jpg"
